$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.089251970589998564
$ws.Range("A2").Value = -0.0099999998471282936
$ws.Range("A3").Value = -0.0089999998482390708
$ws.Range("A4").Value = 0.28399773688155605
$ws.Range("A5").Value = -0.0059999998520590125
$ws.Range("A6").Value = -0.0059999998467077376
$ws.Range("A7").Value = -0.019999999820512926
$ws.Range("A8").Value = -0.019999999818776537
$ws.Range("A9").Value = -0.0059999998422730627
$ws.Range("A10").Value = -0.0059999998404478561
$ws.Range("A11").Value = 0.0045821768697855703
$ws.Range("A12").Value = -0.0059999998398532206
$ws.Range("A13").Value = -0.0059999998376243369
$ws.Range("A14").Value = -0.011999999825985874
$ws.Range("A15").Value = -0.0038239362800052135
$ws.Range("A16").Value = -0.0059999998359669959
$ws.Range("A17").Value = -0.0059999998353124084
$ws.Range("A18").Value = -0.0089999998296708128
$ws.Range("A19").Value = -0.0089999998489185273
$ws.Range("A20").Value = -0.0089999998476599785
$ws.Range("A21").Value = -0.0089999998474965537
$ws.Range("A22").Value = -0.008999999847353557
$ws.Range("A23").Value = -0.032443850751516123
$ws.Range("A24").Value = -0.04199999978403568
$ws.Range("A25").Value = -0.041999999782953878
$ws.Range("A26").Value = -0.0059999998462672011
$ws.Range("A27").Value = -0.0059999998455628756
$ws.Range("A28").Value = -0.0059999998428112988
$ws.Range("A29").Value = -0.011999999829988894
$ws.Range("A30").Value = -0.019999999814461766
$ws.Range("A31").Value = -0.01499999982241107
$ws.Range("A32").Value = -0.014762045723988315
$ws.Range("A33").Value = -0.0059999998386244258
